# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (holdings detail, same layout as the
# other quarterly sheets) right before the "总计" (totals) sheet, and adds
# a matching summary row at the top of "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet, positioned right before "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$styleHdr = $totalSheet.Range("B1").Style
$styleIdx = $totalSheet.Range("A2").Style

$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$q1.Range("B1:H1").Style = $styleHdr

# Data rows: the figures in columns D-G are kept as text, matching how
# this workbook stores them throughout (formatted percentages/amounts).
$q1.Range("A2").Style = $styleIdx
$q1.Range("A2").Value = 0
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value = "160726"
$q1.Range("C2").Value = "嘉实瑞享定期开放灵活配置混合"
$q1.Range("D2").Value = "23.58"
$q1.Range("E2").Value = "63.95"
$q1.Range("F2").Value = "2.28"
$q1.Range("G2").Value = "0.5376"
$q1.Range("H2").Value = 8

$q1.Range("A3").Style = $styleIdx
$q1.Range("A3").Value = 1
$q1.Range("B3:G3").NumberFormat = "@"
$q1.Range("B3").Value = "001900"
$q1.Range("C3").Value = "诺安精选价值混合"
$q1.Range("D3").Value = "0.13"
$q1.Range("E3").Value = "89.96"
$q1.Range("F3").Value = "3.66"
$q1.Range("G3").Value = "0.0048"
$q1.Range("H3").Value = 5

# ---------------------------------------------------------------------
# 2) Prepend a 2022-Q1 summary row to "总计", shifting the rest down
# ---------------------------------------------------------------------
# Re-fetch by name: after the insert above, the old $totalSheet reference
# now tracks whatever sheet occupies its original tab position.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Style = $styleIdx
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.54

# Renumber the 0-based index column for the rows pushed down by the insert
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
